$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8: Magnesium chloride unit price - lower bound value updated, G/I become plain values
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Row 9: Zinc sulfate unit price - lower bound value updated, G/I become plain values
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# Row 17/18: Fermentation magnesium/zinc sulfate loading - G formula multiplier changed from 0.08 to 0.8
$ws.Range("G17").Formula = "=E17*0.8"
$ws.Range("G18").Formula = "=E18*0.8"

# Update selection to reflect rows 17:18 being the last edited area
$ws.Range("A17:XFD18").Select()
$ws.Range("A18").Activate()
